# Workbook was edited: update evaluation results on the QuantitativeMetrics sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Execution metrics section
$ws.Range("B5").Value = "no"
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Syntax & Semantic similarity section - updated Code BLEU evaluation
$ws.Range("B12").Value = 0.2891638084736728
$ws.Range("C12").Value = "{'codebleu': 0.28916380847367285, 'ngram_match_score': 0.08533426372796167, 'weighted_ngram_match_score': 0.09334208630549293, 'syntax_match_score': 0.5897435897435898, 'dataflow_match_score': 0.38823529411764707}"

# Update the active selection shown when the sheet is opened
[void]$ws.Range("B6").Select()
